$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-13 (Sending cluster, Ligand, Receptor, Target cluster, then 16 metric columns E:T)
$data = @(
    @("ECs", "Efna4", "Epha7", "ECs", 3, 1, 0.828516, 2.485548, 0.4625620436231038, 0.4821955800271095, 1, 0.3333333333333333, 0.014513, 0.043539, 0.006538124153480057, 0.009244530017173054, 0.012024252708, 0.108218274372, 0.00302428806989531, 0.004457671513708785),
    @("ECs", "Efna4", "Epha7", "FAPs", 3, 1, 0.828516, 2.485548, 0.4625620436231038, 0.4821955800271095, 3, 1, 0.2421403333333333, 0.726421, 0.1090845147039467, 0.154239204841748, 0.200617140412, 1.805554263708, 0.05045835604909211, 0.07437346284158684),
    @("ECs", "Efna4", "Epha7", "M1", 3, 1, 0.828516, 2.485548, 0.4625620436231038, 0.4821955800271095, 1, 0.3333333333333333, 0.013551, 0.040653, 0.006104741983312083, 0.008631752653669954, 0.011227220316, 0.101044982844, 0.002823821927592597, 0.004162192977486926),
    @("ECs", "Efna4", "Epha7", "sCs", 3, 1, 0.828516, 2.485548, 0.4625620436231038, 0.4821955800271095, 2, 1, 1.9495455, 3.899091, 0.878272619159261, 0.827884512487409, 1.615229639478, 9.691377836868002, 0.4062555775765237, 0.399202252694327),
    @("FAPs", "Efna4", "Epha7", "ECs", 3, 1, 0.7438396666666667, 2.231519, 0.4152870872032183, 0.4329140288365043, 1, 0.3333333333333333, 0.014513, 0.043539, 0.006538124153480057, 0.009244530017173054, 0.01079534508233333, 0.097158105741, 0.002715198535471741, 0.004002086734434385),
    @("FAPs", "Efna4", "Epha7", "FAPs", 3, 1, 0.7438396666666667, 2.231519, 0.4152870872032183, 0.4329140288365043, 3, 1, 0.2421403333333333, 0.726421, 0.1090845147039467, 0.154239204841748, 0.1801135848332222, 1.621022263499, 0.04530139037037868, 0.06677231557257998),
    @("FAPs", "Efna4", "Epha7", "M1", 3, 1, 0.7438396666666667, 2.231519, 0.4152870872032183, 0.4329140288365043, 1, 0.3333333333333333, 0.013551, 0.040653, 0.006104741983312083, 0.008631752653669954, 0.010079771323, 0.090717941907, 0.002535220516376873, 0.003736806817220447),
    @("FAPs", "Efna4", "Epha7", "sCs", 3, 1, 0.7438396666666667, 2.231519, 0.4152870872032183, 0.4329140288365043, 2, 1, 1.9495455, 3.899091, 0.878272619159261, 0.827884512487409, 1.4501492748715, 8.700895649229, 0.364735277780991, 0.3584028197122694),
    @("sCs", "Efna4", "Epha7", "ECs", 2, 1, 0.21879, 0.43758, 0.1221508691736778, 0.08489039113638626, 1, 0.3333333333333333, 0.014513, 0.043539, 0.006538124153480057, 0.009244530017173054, 0.00317529927, 0.01905179562, 0.0007986375481130057, 0.0007847717690298841),
    @("sCs", "Efna4", "Epha7", "FAPs", 2, 1, 0.21879, 0.43758, 0.1221508691736778, 0.08489039113638626, 3, 1, 0.2421403333333333, 0.726421, 0.1090845147039467, 0.154239204841748, 0.05297788352999999, 0.31786730118, 0.01332476828447593, 0.01309342642758119),
    @("sCs", "Efna4", "Epha7", "M1", 2, 1, 0.21879, 0.43758, 0.1221508691736778, 0.08489039113638626, 1, 0.3333333333333333, 0.013551, 0.040653, 0.006104741983312083, 0.008631752653669954, 0.00296482329, 0.01778893974, 0.0007456995393426128, 0.0007327528589625825),
    @("sCs", "Efna4", "Epha7", "sCs", 2, 1, 0.21879, 0.43758, 0.1221508691736778, 0.08489039113638626, 2, 1, 1.9495455, 3.899091, 0.878272619159261, 0.827884512487409, 0.426541059945, 1.70616423978, 0.1072817638017463, 0.0702794400808126)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    for ($j = 0; $j -lt 16; $j++) {
        $ws.Cells.Item($r, 5 + $j).Value = $row[4 + $j]
    }
}

Write-Output "rows written"